$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the training/test column (F) for the rows that were still blank
$values = @{
    3  = "training"
    4  = "test"
    5  = "test"
    7  = "training"
    8  = "training"
    9  = "training"
    10 = "test"
    11 = "test"
    14 = "test"
    15 = "test"
    16 = "test"
    17 = "training"
    27 = "test"
    28 = "training"
    29 = "training"
    30 = "training"
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}

# Select C17 to match the recorded selection state
$ws.Range("C17").Select()

# Turn on AutoFilter over the used range
$ws.Range("A1:F37").AutoFilter() | Out-Null

# Register the (hidden, sheet-scoped) _FilterDatabase defined name that Excel
# writes out alongside an AutoFilter
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:F37"))
$filterName.Visible = $false
